$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 (California - Los Angeles) updated figures
$ws.Range("B17").Value = 44033
$ws.Range("C17").Value = 164870
$ws.Range("D17").Value = 4213
$ws.Range("E17").Value = 4391
$ws.Range("F17").Value = 419
$ws.Range("H17").Value = 10.67
$ws.Range("K17").Value = 94674
$ws.Range("L17").Value = 3927

# Row 39 (Delaware) status message updated
$ws.Range("O39").Value = "An error occurred. ... HTTPError('504 Server Error: Gateway Time-out for url: https://myhealthycommunity.dhss.delaware.gov/locations/state/')"
